$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows for Lenah Cheloti, Moses Ngugi and Ochieng Charles are being
# reordered: Ochieng Charles's data moves up to row 3, Lenah Cheloti's data
# moves to row 4, and Moses Ngugi's data moves to row 5.
#
# Numeric-looking values are stored as plain text in this sheet, so each
# such cell is briefly switched to Text format before being assigned (this
# stops Excel from auto-converting "6.00" / "17.14%" into real numbers),
# then the format is cleared again so no stray number-format style is left
# behind on the cell.

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row 3 <- Ochieng Charles
$ws.Cells.Item(3, 1).Value = "Ochieng Charles"
Set-TextValue $ws.Cells.Item(3, 2) "6.00"
Set-TextValue $ws.Cells.Item(3, 3) "35.00"
Set-TextValue $ws.Cells.Item(3, 4) "-29.00"
Set-TextValue $ws.Cells.Item(3, 5) "17.14%"

# Row 4 <- Lenah Cheloti
$ws.Cells.Item(4, 1).Value = "Lenah Cheloti"
Set-TextValue $ws.Cells.Item(4, 2) "3.00"
Set-TextValue $ws.Cells.Item(4, 3) "28.00"
Set-TextValue $ws.Cells.Item(4, 4) "-25.00"
Set-TextValue $ws.Cells.Item(4, 5) "10.71%"

# Row 5 <- Moses  Ngugi
$ws.Cells.Item(5, 1).Value = "Moses  Ngugi"
Set-TextValue $ws.Cells.Item(5, 2) "2.00"
Set-TextValue $ws.Cells.Item(5, 3) "30.00"
Set-TextValue $ws.Cells.Item(5, 4) "-28.00"
Set-TextValue $ws.Cells.Item(5, 5) "6.67%"
